# Added tests for IT6
# - Renumber the step indices for the (previously skipped) step 6 rows
# - Fill in the CookController ("S") test results for steps 4 and 5
# - Widen the newly meaningful "Output" comment column (N) for the added notes
# - Update the active selection to reflect where the edit was made
#
# (Note: the source revision also shows every row's height/descent ratio
# shifting by a uniform ~1.0417x and the workbook's absPath/revision GUID
# changing. Those are artifacts of the file being resaved by a different
# Excel install/user, not of any cell edit, so they are not reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Fill in the missing CookController test marks ("S") for step 4 (row 8)
# and step 5 (row 9), in column L.
$ws.Range("L8").Value = "S"
$ws.Range("L9").Value = "S"

# Renumber steps 7/8/9 down to 6/7/8 now that step 6 has its own tests.
$ws.Range("A10").Value = 6
$ws.Range("A11").Value = 7
$ws.Range("A12").Value = 8

# Widen column N (used for free-text notes about the tests) so the longer
# note added alongside this change is readable.
$ws.Columns.Item(14).ColumnWidth = 62.666666666666664

# Reflect the cell that was being worked on when the file was saved.
$ws.Range("N12").Select()
